$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("B2").Value = 79243
$ws.Range("B3").Value = 79243
$ws.Range("B4").Value = 79243
$ws.Range("B5").Value = 91828
$ws.Range("B6").Value = 79243
$ws.Range("B7").Value = 80348
$ws.Range("B8").Value = 79714
$ws.Range("B9").Value = 58520
$ws.Range("B10").Value = 79244
$ws.Range("B11").Value = 79244
$ws.Range("B12").Value = 80348
